$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 1.214185333333333
$ws.Cells.Item(2, 8).Value = 3.642556
$ws.Cells.Item(2, 9).Value = 0.02886664532879162
$ws.Cells.Item(2, 10).Value = 0.02886664532879162
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 8.309350333333333
$ws.Cells.Item(2, 14).Value = 24.928051
$ws.Cells.Item(2, 15).Value = 0.1535033474258946
$ws.Cells.Item(2, 16).Value = 0.1535033474258946
$ws.Cells.Item(2, 17).Value = 10.08909130426178
$ws.Cells.Item(2, 18).Value = 90.801821738356
$ws.Cells.Item(2, 19).Value = 0.004431126686925578
$ws.Cells.Item(2, 20).Value = 0.004431126686925578

$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 1.214185333333333
$ws.Cells.Item(3, 8).Value = 3.642556
$ws.Cells.Item(3, 9).Value = 0.02886664532879162
$ws.Cells.Item(3, 10).Value = 0.02886664532879162
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 37.153391
$ws.Cells.Item(3, 14).Value = 111.460173
$ws.Cells.Item(3, 15).Value = 0.6863556906301786
$ws.Cells.Item(3, 16).Value = 0.6863556906301786
$ws.Cells.Item(3, 17).Value = 45.11110243579866
$ws.Cells.Item(3, 18).Value = 405.999921922188
$ws.Cells.Item(3, 19).Value = 0.01981278629081919
$ws.Cells.Item(3, 20).Value = 0.01981278629081919

$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 1.214185333333333
$ws.Cells.Item(4, 8).Value = 3.642556
$ws.Cells.Item(4, 9).Value = 0.02886664532879162
$ws.Cells.Item(4, 10).Value = 0.02886664532879162
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 8.668653666666666
$ws.Cells.Item(4, 14).Value = 26.005961
$ws.Cells.Item(4, 15).Value = 0.1601409619439267
$ws.Cells.Item(4, 16).Value = 0.1601409619439267
$ws.Cells.Item(4, 17).Value = 10.52535214181289
$ws.Cells.Item(4, 18).Value = 94.728169276316
$ws.Cells.Item(4, 19).Value = 0.00462273235104685
$ws.Cells.Item(4, 20).Value = 0.00462273235104685

$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 28.525713
$ws.Cells.Item(5, 8).Value = 85.577139
$ws.Cells.Item(5, 9).Value = 0.6781844725971822
$ws.Cells.Item(5, 10).Value = 0.6781844725971822
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 8.309350333333333
$ws.Cells.Item(5, 14).Value = 24.928051
$ws.Cells.Item(5, 15).Value = 0.1535033474258946
$ws.Cells.Item(5, 16).Value = 0.1535033474258946
$ws.Cells.Item(5, 17).Value = 237.030142825121
$ws.Cells.Item(5, 18).Value = 2133.271285426089
$ws.Cells.Item(5, 19).Value = 0.1041035867159323
$ws.Cells.Item(5, 20).Value = 0.1041035867159324

$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 28.525713
$ws.Cells.Item(6, 8).Value = 85.577139
$ws.Cells.Item(6, 9).Value = 0.6781844725971822
$ws.Cells.Item(6, 10).Value = 0.6781844725971822
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 37.153391
$ws.Cells.Item(6, 14).Value = 111.460173
$ws.Cells.Item(6, 15).Value = 0.6863556906301786
$ws.Cells.Item(6, 16).Value = 0.6863556906301786
$ws.Cells.Item(6, 17).Value = 1059.826968642783
$ws.Cells.Item(6, 18).Value = 9538.442717785047
$ws.Cells.Item(6, 19).Value = 0.4654757720641025
$ws.Cells.Item(6, 20).Value = 0.4654757720641025

$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 28.525713
$ws.Cells.Item(7, 8).Value = 85.577139
$ws.Cells.Item(7, 9).Value = 0.6781844725971822
$ws.Cells.Item(7, 10).Value = 0.6781844725971822
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 8.668653666666666
$ws.Cells.Item(7, 14).Value = 26.005961
$ws.Cells.Item(7, 15).Value = 0.1601409619439267
$ws.Cells.Item(7, 16).Value = 0.1601409619439267
$ws.Cells.Item(7, 17).Value = 247.279526591731
$ws.Cells.Item(7, 18).Value = 2225.515739325579
$ws.Cells.Item(7, 19).Value = 0.1086051138171474
$ws.Cells.Item(7, 20).Value = 0.1086051138171474

$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 12.32198033333333
$ws.Cells.Item(8, 8).Value = 36.965941
$ws.Cells.Item(8, 9).Value = 0.2929488820740263
$ws.Cells.Item(8, 10).Value = 0.2929488820740263
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 8.309350333333333
$ws.Cells.Item(8, 14).Value = 24.928051
$ws.Cells.Item(8, 15).Value = 0.1535033474258946
$ws.Cells.Item(8, 16).Value = 0.1535033474258946
$ws.Cells.Item(8, 17).Value = 102.3876513901101
$ws.Cells.Item(8, 18).Value = 921.488862510991
$ws.Cells.Item(8, 19).Value = 0.04496863402303667
$ws.Cells.Item(8, 20).Value = 0.04496863402303668

$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 12.32198033333333
$ws.Cells.Item(9, 8).Value = 36.965941
$ws.Cells.Item(9, 9).Value = 0.2929488820740263
$ws.Cells.Item(9, 10).Value = 0.2929488820740263
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 37.153391
$ws.Cells.Item(9, 14).Value = 111.460173
$ws.Cells.Item(9, 15).Value = 0.6863556906301786
$ws.Cells.Item(9, 16).Value = 0.6863556906301786
$ws.Cells.Item(9, 17).Value = 457.8033532186437
$ws.Cells.Item(9, 18).Value = 4120.230178967793
$ws.Cells.Item(9, 19).Value = 0.201067132275257
$ws.Cells.Item(9, 20).Value = 0.201067132275257

$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 12.32198033333333
$ws.Cells.Item(10, 8).Value = 36.965941
$ws.Cells.Item(10, 9).Value = 0.2929488820740263
$ws.Cells.Item(10, 10).Value = 0.2929488820740263
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 8.668653666666666
$ws.Cells.Item(10, 14).Value = 26.005961
$ws.Cells.Item(10, 15).Value = 0.1601409619439267
$ws.Cells.Item(10, 16).Value = 0.1601409619439267
$ws.Cells.Item(10, 17).Value = 106.8149799971445
$ws.Cells.Item(10, 18).Value = 961.334819974301
$ws.Cells.Item(10, 19).Value = 0.04691311577573252
$ws.Cells.Item(10, 20).Value = 0.04691311577573252
